$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.330.52'
$ws.Range("E2").Value = '  +4.23%  '
$ws.Range("D3").Value = '3.053.15'
$ws.Range("E3").Value = '  +2.47%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''548.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.30%  '
$ws.Range("D6").Value = '''139.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.60%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.049.35'
$ws.Range("E8").Value = '  +2.40%  '
$ws.Range("D9").Value = '''0.499'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.98%  '
$ws.Range("D10").Value = '''6.26'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.50%  '
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("E12").Value = '  +2.77%  '
$ws.Range("D13").Value = '''0.0000226'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.71%  '
$ws.Range("D14").Value = '''34.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.31%  '
$ws.Range("D15").Value = '3.557.48'
$ws.Range("E15").Value = '  +3.13%  '
$ws.Range("D16").Value = '63.389.88'
$ws.Range("E16").Value = '  +4.35%  '
$ws.Range("D17").Value = '3.058.78'
$ws.Range("E17").Value = '  +2.87%  '
$ws.Range("E18").Value = '  -1.19%  '
$ws.Range("D19").Value = '''6.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.27%  '
$ws.Range("D20").Value = '''478.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.64%  '
$ws.Range("D21").Value = '''13.51'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.57%  '
$ws.Range("D22").Value = '''0.672'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("D23").Value = '''7.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.84%  '
$ws.Range("D24").Value = '''80.83'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.46%  '
$ws.Range("D25").Value = '''12.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.20%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("D27").Value = '''2.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.39%  '
$ws.Range("D28").Value = '''7.89'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.51%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("D30").Value = '''1.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.50%  '
$ws.Range("D31").Value = '''25.90'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.98%  '
$ws.Range("E32").Value = '  +0.83%  '
$ws.Range("D33").Value = '''2.40'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.28%  '
$ws.Range("D34").Value = '''5.63'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.76%  '
$ws.Range("D35").Value = '''55.67'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.51%  '
$ws.Range("E36").Value = '  +3.79%  '
$ws.Range("D37").Value = '''461.47'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.82%  '
$ws.Range("D38").Value = '''0.0810'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.83%  '
$ws.Range("D39").Value = '3.117.33'
$ws.Range("E39").Value = '  -1.90%  '
$ws.Range("D40").Value = '''0.0392'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.05%  '
$ws.Range("E41").Value = '  +1.69%  '
$ws.Range("D42").Value = '''8.20'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.45%  '
$ws.Range("E43").Value = '  +7.11%  '
$ws.Range("D44").Value = '''27.90'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +10.38%  '
$ws.Range("D45").Value = '''0.250'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.44%  '
$ws.Range("E47").Value = '  +4.95%  '
$ws.Range("E48").Value = '  +1.24%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = '''115.67'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.93%  '
$ws.Range("B50").Value = 'PEPE'
$ws.Range("C50").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D50").Value = '0.0₃0506'
$ws.Range("E50").Value = '  +2.60%  '
$ws.Range("E51").Value = '  +5.19%  '
